$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary rows 10-12: pick up the style already used for the row-9/row-15
# section headers ("mtitleStyle", cellXfs index 4) for the label cells in
# column A, and refresh the numeric grading figures now that the grader
# copes with float input.
# ---------------------------------------------------------------------------
$ws.Range("A15").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A10").Value = "No."
$ws.Range("A11").Value = "Marking"
$ws.Range("A12").Value = "Total"

$ws.Range("B10").Value = 25
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100/112"

# ---------------------------------------------------------------------------
# The per-question breakdown used to be split across three side-by-side
# blocks (A:B, D:E, G:H). It now only needs one block (A:B) plus a couple of
# leftover rows in D:E, so drop the G:H block entirely and the now-unused
# D:E rows below it.
# ---------------------------------------------------------------------------
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

$ws.Range("B11").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A20").Value = "Option B"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("A24").Value = "Option A"
$ws.Range("A25").Value = "Option A"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option B"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option A"
$ws.Range("A37").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"
